$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.46327006816864
$ws.Range("B1").Value = 3.094699859619141
$ws.Range("C1").Value = 1.743524074554443
$ws.Range("D1").Value = 1.289110064506531
$ws.Range("E1").Value = 1.122492790222168
